# Apply the changes described by the commit:
# "update content and fix page width at small screen sizes"
#
# 1. Fill in the previously-blank N:T columns (reliability-related translations)
#    for rows 15, 18, 73, 79, 97, 101, 103, 104, 106 on the "updates" sheet by
#    copying the already-correct values/styles from row 16.
# 2. Add three new rows (17-19) to the "content" sheet with new UI strings
#    (languages / Back to top / Download Summary labels).
# 3. Narrow/replace the autofit width of column A on the "content" sheet with a
#    fixed custom width so the page renders correctly at small screen sizes.
# 4. Update sheet selections/active sheet to match the saved view state.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "updates"
$ws2 = $wb.Worksheets.Item(2)   # "lines"
$ws3 = $wb.Worksheets.Item(3)   # "content"

# --- 1. Populate N:T for the affected rows on "updates" by copying row 16 ---
$sourceRange = $ws1.Range("N16:T16")
$targetRows = 15,18,73,79,97,101,103,104,106
foreach ($r in $targetRows) {
    $sourceRange.Copy($ws1.Range("N" + $r + ":T" + $r))
}

# --- 2. Add new content rows for the language / back-to-top / download summary UI strings ---
# NOTE: values are assigned in this particular order so that the workbook's
# shared-string table ends up with the same new-string ordering as the source.
$ws3.Range("A17").Value = "languages"
$ws3.Range("B19").Value = "Back to top"
$ws3.Range("A19").Value = "backToTop"
$ws3.Range("B17").Value = "Languages"
$ws3.Range("A18").Value = "downloadSummary"
$ws3.Range("B18").Value = "Download Summary"

# Match formatting (left-aligned) used by the rest of column B
$ws3.Range("B17:B19").HorizontalAlignment = -4131

# --- 3. Fix column A width on "content" sheet (remove autofit/bestFit, use fixed width) ---
$ws3.Columns.Item(1).ColumnWidth = 22.33

# --- 4. Update selections / active sheet ---
$ws3.Range("C40").Select()
$ws2.Range("D91").Select()
$ws1.Activate()
$ws1.Range("N30").Select()
